# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 3 (pushing the existing
# rows 3..63 down to 4..64) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3; this shifts rows 3..63 down
# to 4..64 (preserving their formatting/styles, e.g. the date format on
# column D) and grows the used range to A1:T64 automatically.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's new record.
$ws.Cells.Item(3, 1).Value  = 11
$ws.Cells.Item(3, 2).Value  = 'Vega Monumental Concepción'
$ws.Cells.Item(3, 3).Value  = 'Bíobío'
$ws.Cells.Item(3, 4).Value  = (Get-Date -Year 2023 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(3, 5).Value  = 8
$ws.Cells.Item(3, 6).Value  = 'Fruta'
$ws.Cells.Item(3, 7).Value  = 100107
$ws.Cells.Item(3, 8).Value  = 'Otros'
$ws.Cells.Item(3, 9).Value  = 100107002
$ws.Cells.Item(3, 10).Value = 'Chirimoya'
$ws.Cells.Item(3, 11).Value = 'Cultivar IV Región'
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 19000
$ws.Cells.Item(3, 15).Value = 20000
$ws.Cells.Item(3, 16).Value = 19500
$ws.Cells.Item(3, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(3, 19).Value = 1950
$ws.Cells.Item(3, 20).Value = 10

Write-Host "Inserted new row 3; sheet now spans $($ws.UsedRange.Address())"
